# Auto-update draw results: append the 2025-10-09 Pick 3 row (row 23).
#
# All five columns in this sheet are stored as literal TEXT (t="str"),
# even values that look like dates ("2025-10-09") or numbers ("251009").
# Assigning those strings straight to Range.Value would make Excel's COM
# layer auto-coerce them into a date serial / a number, so instead each
# value is written as a `="literal text"` text-formula and then flattened
# back to a plain value via Copy + PasteSpecial(xlPasteValues). That keeps
# the stored cell type as Text without touching any cell's number format
# (NumberFormat stays "General", exactly like the existing rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 23

$ws.Cells.Item($newRow, 1).Formula = "=""2025-10-09"""
$ws.Cells.Item($newRow, 2).Formula = "=""Pick 3"""
$ws.Cells.Item($newRow, 3).Formula = "=""251009"""
$ws.Cells.Item($newRow, 4).Formula = "=""5-9-4"""
$ws.Cells.Item($newRow, 5).Formula = "=""2025-10-09T21:39:03.628+04:00"""

$newRowRange = $ws.Range("A" + $newRow + ":E" + $newRow)
$newRowRange.Copy()
$newRowRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# The source file also widens its "numbers stored as text, please don't
# warn" ignored range from A1:E22 to A1:E23 so the newly appended row
# doesn't show the green-triangle warning either. Mirror that intent here
# (xlNumberAsText = -2146826246); harmless if the host doesn't persist it.
try {
    $fullRange = $ws.Range("A1:E" + $newRow)
    $fullRange.Errors.Item(-2146826246).Ignore = $true
} catch {
}
